# CoursesDT, status and notes
# Update the "Status" column (AA) so the default status text is
# "Not processed" while the rows that finished processing successfully
# keep showing "Success".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that succeeded (CHEE2501 TUT1/LEC1/FLD1/EXC1, ECON1120 LEC1) retain "Success"
$successRows = @(11, 12, 15, 16, 17, 18)

# All data rows (2-21) get a Status value; default is now "Not processed"
for ($r = 2; $r -le 21; $r++) {
    if ($successRows -contains $r) {
        $ws.Range("AA$r").Value = "Success"
    } else {
        $ws.Range("AA$r").Value = "Not processed"
    }
}
